$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(7, 10, 9, 14, 15),
    @(10, 19, 22, 13, 23),
    @(8, 3, 16, 19, 21),
    @(19, 12, 5, 23, 4),
    @(16, 22, 18, 14, 20)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
